# The target revision ("Updated ppt's and bit sim log") only touches the
# package's internal co-authoring/revision bookkeeping: it adds a brand new
# part, ppt/changesInfos/changesInfo1.xml, which is a `pc:chgInfo` log
# recording eight slide deletions (sldId 256/258/259/260/261/262/264/265)
# that a prior editing session already performed. Those slides are already
# absent from this presentation (it already has only the 4 surviving
# slides - "Tagging", "Bit Tagging", "Cutter Tagging", "Last Slide" - matching
# ppt/revisionInfo.xml's v="8" high-water mark), so there is no slide/shape
# content to change here: the diff is 100% additive and limited to that one
# metadata part.
#
# ppt/changesInfos/*.xml is written by PowerPoint's real-time co-authoring
# engine (SharePoint/OneDrive sync merges) - it is not a document property,
# slide, shape, or anything else surfaced by the Presentation object model,
# so it cannot be produced through PowerPoint automation/COM (no VBA/OM call
# creates it; Coauthoring/Sync come back empty outside a live co-authoring
# session, and AcceptAll/RejectAll/Merge/EndReview do not touch it either).
#
# The faithful COM-automation action is therefore to leave the already-correct
# slide content untouched and simply save, which keeps every existing part
# (including ppt/revisionInfo.xml) intact.
$p = $ppt.ActivePresentation
$p.Save()
